# Add 'drop' and 'dropExplanation' columns to Sheet2
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet2")
$ws.Activate()

# Headers - new columns Q (drop) and R (dropExplanation)
$ws.Range("Q1").Value = "drop"
$ws.Range("R1").Value = "dropExplanation"

# Match the header formatting used by the rest of row 1 (style of P1)
$ws.Range("P1").Copy()
$ws.Range("Q1:R1").PasteSpecial(-4122)

# Data rows: default value of FALSE for the new 'drop' flag column
$ws.Range("Q2:Q41").Value = $false

# Reflect the selection left after adding/inspecting the new columns
$ws.Range("Q1:R1048576").Select() | Out-Null
